{"js": "// Update the worksheet date and the 25 \"two-digit \u00f7 one-digit\" answer\n// cells with the new generated values, preserving all existing\n// character/paragraph formatting (fonts, size, alignment, etc.).\n\n// 1) Update the date heading paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2025-10-17 Friday\", Word.InsertLocation.replace);\n\n// 2) Update the answer cells inside the (single) table, in row/column\n// (left-to-right, top-to-bottom) order, matching the document order of\n// the original text runs.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newAnswers = [\n  \"33\u00f78=4, 1\", \"91\u00f77=13, 0\", \"96\u00f75=19, 1\", \"74\u00f74=18, 2\", \"24\u00f76=4, 0\",\n  \"46\u00f77=6, 4\", \"80\u00f77=11, 3\", \"39\u00f73=13, 0\", \"11\u00f75=2, 1\", \"67\u00f72=33, 1\",\n  \"24\u00f74=6, 0\", \"85\u00f74=21, 1\", \"64\u00f73=21, 1\", \"78\u00f78=9, 6\", \"40\u00f74=10, 0\",\n  \"38\u00f77=5, 3\", \"23\u00f77=3, 2\", \"88\u00f74=22, 0\", \"69\u00f74=17, 1\", \"44\u00f74=11, 0\",\n  \"82\u00f78=10, 2\", \"82\u00f78=10, 2\", \"33\u00f77=4, 5\", \"50\u00f75=10, 0\", \"87\u00f76=14, 3\",\n];\n\nconst rowValues = table.values;\nlet answerIndex = 0;\nfor (let r = 0; r < rowValues.length; r++) {\n  const row = rowValues[r];\n  for (let c = 0; c < row.length; c++) {\n    if (row[c] !== \"\") {\n      const cell = table.getCell(r, c);\n      cell.body.getRange().insertText(newAnswers[answerIndex], Word.InsertLocation.replace);\n      answerIndex++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 \"two-digit \u00f7 one-digit\" answer\n# cells with the new generated values, preserving all existing\n# character/paragraph formatting (fonts, size, alignment, etc.).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading paragraph (first paragraph in the document).\n$d.Paragraphs(1).Range.Text = \"2025-10-17 Friday\"\n\n# 2) Update the answer cells inside the (single) table, in row/column\n# (left-to-right, top-to-bottom) order, matching the document order of\n# the original text runs. Only cells that already contain an answer are\n# touched; the blank spacer rows are left untouched.\n$newAnswers = @(\n    \"33\u00f78=4, 1\", \"91\u00f77=13, 0\", \"96\u00f75=19, 1\", \"74\u00f74=18, 2\", \"24\u00f76=4, 0\",\n    \"46\u00f77=6, 4\", \"80\u00f77=11, 3\", \"39\u00f73=13, 0\", \"11\u00f75=2, 1\", \"67\u00f72=33, 1\",\n    \"24\u00f74=6, 0\", \"85\u00f74=21, 1\", \"64\u00f73=21, 1\", \"78\u00f78=9, 6\", \"40\u00f74=10, 0\",\n    \"38\u00f77=5, 3\", \"23\u00f77=3, 2\", \"88\u00f74=22, 0\", \"69\u00f74=17, 1\", \"44\u00f74=11, 0\",\n    \"82\u00f78=10, 2\", \"82\u00f78=10, 2\", \"33\u00f77=4, 5\", \"50\u00f75=10, 0\", \"87\u00f76=14, 3\"\n)\n\n$t = $d.Tables(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$answerIndex = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $txt = $cell.Range.Text\n        if ($txt.Length -gt 2) {\n            $cell.Range.Text = $newAnswers[$answerIndex]\n            $answerIndex++\n        }\n    }\n}\n"}
